# feat: add 2022-Q1 data
#
# The workbook currently ends with a "总计" (grand-total) summary sheet.
# This change:
#   1. Turns that existing sheet into the new "2022-Q1" per-fund holdings
#      sheet (new header row + 4 fund rows).
#   2. Adds a brand new "总计" sheet right after it, containing the old
#      summary rows plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# A sheet that already has the "index column + bold/bordered header" look
# we need to replicate (2021-Q4 holds the 8-row/8-col fund table).
$formatSrc = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet into the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Pull over the header/index-column formatting (bold, centered, thin
# border for row 1 and the A-column row numbers) without minting new
# style entries.
$formatSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$formatSrc.Range("A2:A5").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (0-based row counter), same convention as every other
# per-quarter sheet in this workbook.
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3

# Fund rows. 基金代码/基金规模/股票总仓位/仓位占比/持有市值(亿元) are all
# stored as text (fund codes need their leading zeros kept, the numeric
# ones keep trailing zeros like "14.75"), so force text entry with a
# leading apostrophe the same way a user typing them in would; 仓位排名
# stays a real number.
$q1.Range("B2").Value = "'006679"
$q1.Range("C2").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A"
$q1.Range("D2").Value = "'14.75"
$q1.Range("E2").Value = "'83.19"
$q1.Range("F2").Value = "'5.87"
$q1.Range("G2").Value = "'0.8658"
$q1.Range("H2").Value = 4

$q1.Range("B3").Value = "'162719"
$q1.Range("C3").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
$q1.Range("D3").Value = "'14.75"
$q1.Range("E3").Value = "'83.19"
$q1.Range("F3").Value = "'5.87"
$q1.Range("G3").Value = "'0.8658"
$q1.Range("H3").Value = 4

$q1.Range("B4").Value = "'006680"
$q1.Range("C4").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C"
$q1.Range("D4").Value = "'4.73"
$q1.Range("E4").Value = "'83.19"
$q1.Range("F4").Value = "'5.87"
$q1.Range("G4").Value = "'0.2777"
$q1.Range("H4").Value = 4

$q1.Range("B5").Value = "'004243"
$q1.Range("C5").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
$q1.Range("D5").Value = "'4.73"
$q1.Range("E5").Value = "'83.19"
$q1.Range("F5").Value = "'5.87"
$q1.Range("G5").Value = "'0.2777"
$q1.Range("H5").Value = 4

# ---------------------------------------------------------------------
# Step 2: add the new "总计" sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$formatSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$formatSrc.Range("A2:A7").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 2.29

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 2.33

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.83

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 4
$total.Range("D5").Value = 0.95

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 1.24

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 7
$total.Range("D7").Value = 4.04
